$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
$c1 = $tcs.Colors(3)
Write-Host "Before dk2: $($c1.RGB)"
$c1.RGB = 255
Write-Host "After dk2: $($c1.RGB)"
